# "made stations sizing similar to dashboard"
#
# The workbook has two sheets:
#   1) "Example" - a per-sprint effort log
#   2) "Count"   - a per-sprint tally / totals sheet
#
# The project was renamed from the placeholder "Example Template" /
# "Example Project" to "SkyPi", Sprint #3 ("stations page foundation" work)
# got logged with real data (mirroring how the dashboard/login sprint was
# already sized on the Count sheet), Sprint #4's effort hours were refined
# (8 -> 8.5), and the Count sheet's totals row became live SUM() formulas
# instead of hand-typed numbers. The active tab also moved from "Example"
# over to "Count".

$wb = $excel.ActiveWorkbook
$wsExample = $wb.Worksheets.Item("Example")
$wsCount = $wb.Worksheets.Item("Count")

# ---------------------------------------------------------------------
# "Example" sheet
# ---------------------------------------------------------------------

# Project name: Example Template -> SkyPi
$wsExample.Range("C2").Value = "SkyPi"

# Sprint # 3 (row 6) was logged blank before - fill in the real entry,
# sized/shaped like the dashboard sprint rows above it.
$wsExample.Range("C6").Value = 20
$wsExample.Range("E6").Value = "Worked on more tutorials, helped build dashboard page, helped build login page, meet with sponser"
$wsExample.Range("F6").Value = "dashboard.js`nlogin.js`nprojectmanager"
$wsExample.Range("G6").Value = "skypi\fronted\client\src`nsprint03\projectmanager"
$wsExample.Range("H6").Value = "https://www.youtube.com/watch?v=A71aqufiNtQ"
$wsExample.Rows.Item(6).RowHeight = 78.75

# Sprint # 4 (row 7): effort hours refined from 8 to 8.5
$wsExample.Range("C7").Value = 8.5

# ---------------------------------------------------------------------
# "Count" sheet
# ---------------------------------------------------------------------

# Project name: Example Project -> SkyPi
$wsCount.Range("C2").Value = "SkyPi"

# Totals row (row 5): switch from hard-coded numbers to live sums over
# the per-sprint rows (6:13), matching each column.
$wsCount.Range("C5").Formula = "=SUM(C6:C13)"
$wsCount.Range("D5:P5").Formula = "=SUM(D6:D13)"

# Row 7 (sprint 2): column D corrected from 2 to 1
$wsCount.Range("D7").Value = 1

# Row 8 (sprint 3) was blank before - now carries the same per-category
# tally sizing as the dashboard/login sprint rows.
$wsCount.Range("C8").Value = 1
$wsCount.Range("D8").Value = 1
$wsCount.Range("E8").Value = 1
$wsCount.Range("F8").Value = 1
$wsCount.Range("G8").Value = 0
$wsCount.Range("H8").Value = 1
$wsCount.Range("I8").Value = 1
$wsCount.Range("J8").Value = 1
$wsCount.Range("K8").Value = 1
$wsCount.Range("L8").Value = 1
$wsCount.Range("M8").Value = 1
$wsCount.Range("N8").Value = 0
$wsCount.Range("O8").Value = 0
$wsCount.Range("P8").Value = 0

# Row 9 (sprint 4): two new tallies recorded
$wsCount.Range("G9").Value = 1
$wsCount.Range("P9").Value = 1

# ---------------------------------------------------------------------
# View state: active tab moves from "Example" to "Count"
# ---------------------------------------------------------------------

$wsCount.Activate()
$excel.ActiveWindow.Zoom = 90
$wsCount.Range("C9").Select()
